# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the target OOXML diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '30.553.45'
$ws.Cells.Item(2,5).Value = '  -0.35%  '

$ws.Cells.Item(3,4).Value = '1.876.96'
$ws.Cells.Item(3,5).Value = '  -0.79%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '0.9993'
$ws.Cells.Item(4,5).Value = '  -0.10%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '236.86'
$ws.Cells.Item(5,5).Value = '  -3.05%  '

$ws.Cells.Item(6,5).Value = '  +0.05%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.4873'
$ws.Cells.Item(7,5).Value = '  -1.74%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.2903'
$ws.Cells.Item(8,5).Value = '  -1.98%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06680'
$ws.Cells.Item(9,5).Value = '  -2.04%  '

$ws.Cells.Item(10,4).Value = '1.873.97'
$ws.Cells.Item(10,5).Value = '  -0.94%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '16.67'
$ws.Cells.Item(11,5).Value = '  -2.72%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.07231'
$ws.Cells.Item(12,5).Value = '  -1.14%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '89.75'
$ws.Cells.Item(13,5).Value = '  -1.60%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '5.011'
$ws.Cells.Item(14,5).Value = '  -1.74%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.6560'
$ws.Cells.Item(15,5).Value = '  -2.77%  '

$ws.Cells.Item(16,4).Value = '30.504.87'
$ws.Cells.Item(16,5).Value = '  -0.42%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.000007831'
$ws.Cells.Item(17,5).Value = '  -1.20%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '0.9997'
$ws.Cells.Item(18,5).Value = '  -0.10%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '13.03'
$ws.Cells.Item(19,5).Value = '  -1.81%  '

$ws.Cells.Item(20,4).Value = '2.115.30'
$ws.Cells.Item(20,5).Value = '  -0.94%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '0.9981'
$ws.Cells.Item(21,5).Value = '  -0.43%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '213.60'
$ws.Cells.Item(22,5).Value = '  +17.80%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '4.739'
$ws.Cells.Item(23,5).Value = '  -2.59%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '6.154'
$ws.Cells.Item(24,5).Value = '  +1.57%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '9.391'
$ws.Cells.Item(25,5).Value = '  +0.86%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '156.78'
$ws.Cells.Item(26,5).Value = '  +1.53%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '19.31'
$ws.Cells.Item(27,5).Value = '  +2.44%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.831'
$ws.Cells.Item(28,5).Value = '  -5.09%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '1.413'
$ws.Cells.Item(29,5).Value = '  +1.87%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '4.269'
$ws.Cells.Item(30,5).Value = '  -1.32%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.09061'
$ws.Cells.Item(31,5).Value = '  +1.28%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '3.932'
$ws.Cells.Item(32,5).Value = '  -2.69%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.05128'
$ws.Cells.Item(33,5).Value = '  -1.53%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '0.7279'
$ws.Cells.Item(34,5).Value = '  -1.39%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.081'
$ws.Cells.Item(35,5).Value = '  -4.64%  '

$ws.Cells.Item(36,5).Value = '  +0.40%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.01816'
$ws.Cells.Item(37,5).Value = '  -2.71%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '2.660'
$ws.Cells.Item(38,5).Value = '  -1.38%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '0.9199'
$ws.Cells.Item(39,5).Value = '  -1.47%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '2.049'
$ws.Cells.Item(40,5).Value = '  -5.47%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.4432'
$ws.Cells.Item(41,5).Value = '  +1.62%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '104.68'
$ws.Cells.Item(42,5).Value = '  -1.18%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '5.750'
$ws.Cells.Item(43,5).Value = '  -1.14%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.9957'
$ws.Cells.Item(44,5).Value = '  -0.48%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.1330'
$ws.Cells.Item(45,5).Value = '  -1.55%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '7.358'
$ws.Cells.Item(46,5).Value = '  -3.89%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.4024'
$ws.Cells.Item(47,5).Value = '  +3.23%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.05834'
$ws.Cells.Item(48,5).Value = '  -0.10%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '8.606'
$ws.Cells.Item(49,5).Value = '  +0.70%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '1.416'
$ws.Cells.Item(50,5).Value = '  +2.57%  '

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '33.32'
$ws.Cells.Item(51,5).Value = '  -0.20%  '
